$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cell values while preserving original text storage
# (apostrophe-prefix forces text interpretation for numeric-looking strings,
# then Style is reset to Normal so no stray formatting is introduced).
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '37.410.50'
Set-TextValue 'E2' '  -1.14%  '
Set-TextValue 'D3' '2.048.90'
Set-TextValue 'E3' '  -2.09%  '
Set-TextValue 'D4' '0.998'
Set-TextValue 'E4' '  -0.30%  '
Set-TextValue 'D5' '229.12'
Set-TextValue 'E5' '  -2.02%  '
Set-TextValue 'E6' '  -2.04%  '
Set-TextValue 'E7' '  +0.04%  '
Set-TextValue 'D8' '56.35'
Set-TextValue 'E8' '  -3.82%  '
Set-TextValue 'D9' '0.386'
Set-TextValue 'E9' '  -2.42%  '
Set-TextValue 'D10' '0.0813'
Set-TextValue 'E10' '  +3.76%  '
Set-TextValue 'E11' '  -1.81%  '
Set-TextValue 'D12' '2.348.91'
Set-TextValue 'E12' '  -2.18%  '
Set-TextValue 'D13' '14.55'
Set-TextValue 'E13' '  -4.39%  '
Set-TextValue 'D14' '20.69'
Set-TextValue 'E14' '  -3.18%  '
Set-TextValue 'D15' '0.756'
Set-TextValue 'E15' '  -3.26%  '
Set-TextValue 'E16' '  -2.11%  '
Set-TextValue 'D17' '2.053.93'
Set-TextValue 'E17' '  -1.52%  '
Set-TextValue 'D18' '37.292.00'
Set-TextValue 'E18' '  -1.33%  '
Set-TextValue 'D19' '6.04'
Set-TextValue 'E19' '  -1.79%  '
Set-TextValue 'D20' '69.86'
Set-TextValue 'E20' '  -1.98%  '
Set-TextValue 'D21' '0.0₃0855'
Set-TextValue 'E21' '  +1.94%  '
Set-TextValue 'D22' '226.20'
Set-TextValue 'E22' '  -1.88%  '
Set-TextValue 'D23' '0.999'
Set-TextValue 'E23' '  +0.18%  '
Set-TextValue 'E24' '  -0.58%  '
Set-TextValue 'E25' '  -4.76%  '
Set-TextValue 'D26' '9.56'
Set-TextValue 'E26' '  -2.83%  '
Set-TextValue 'D27' '168.08'
Set-TextValue 'E27' '  -1.91%  '
Set-TextValue 'E28' '  -3.71%  '
Set-TextValue 'E29' '  -0.37%  '
Set-TextValue 'E30' '  -3.10%  '
Set-TextValue 'E31' '  -2.43%  '
Set-TextValue 'D32' '4.53'
Set-TextValue 'E32' '  -4.13%  '
Set-TextValue 'D33' '0.0613'
Set-TextValue 'E33' '  -3.42%  '
Set-TextValue 'D34' '4.56'
Set-TextValue 'E34' '  -2.42%  '
Set-TextValue 'D35' '2.39'
Set-TextValue 'E35' '  -4.94%  '
Set-TextValue 'E36' '  -0.12%  '
Set-TextValue 'D37' '0.999'
Set-TextValue 'E37' '  -0.13%  '
Set-TextValue 'D38' '3.19'
Set-TextValue 'E38' '  -3.97%  '
Set-TextValue 'D39' '5.43'
Set-TextValue 'E39' '  +0.63%  '
Set-TextValue 'B40' 'VeChain'
Set-TextValue 'C40' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D40' '0.0221'
Set-TextValue 'E40' '  -6.05%  '
Set-TextValue 'B41' 'Maker'
Set-TextValue 'C41' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D41' '1.506.32'
Set-TextValue 'E41' '  +4.02%  '
Set-TextValue 'B42' 'InjectiveProtocol'
Set-TextValue 'C42' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D42' '17.06'
Set-TextValue 'E42' '  +2.09%  '
Set-TextValue 'B43' 'HuobiToken'
Set-TextValue 'C43' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D43' '2.88'
Set-TextValue 'E43' '  -1.97%  '
Set-TextValue 'D44' '96.29'
Set-TextValue 'E44' '  -5.63%  '
Set-TextValue 'D45' '0.0937'
Set-TextValue 'E45' '  -3.94%  '
Set-TextValue 'E46' '  -2.99%  '
Set-TextValue 'D47' '1.02'
Set-TextValue 'E47' '  -4.54%  '
Set-TextValue 'D48' '7.16'
Set-TextValue 'E48' '  -1.54%  '
Set-TextValue 'E49' '  -2.12%  '
Set-TextValue 'D50' '2.234.19'
Set-TextValue 'E50' '  -2.19%  '
Set-TextValue 'D51' '3.62'
Set-TextValue 'E51' '  -12.60%  '
